# "Pavel - new user for linking test"
# Adds a new row to the Users sheet for the "Linking_AutoUser" test account,
# used as the default user for the Linking test suite.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# New row goes right after the current last data row (row 51 -> row 52).
$newRow = 52

# Clone formatting (borders/font/number-format) from an existing plain data
# row so the new row visually matches the rest of the table, then overwrite
# the values for the new user.
[void]$ws.Range("A29:G29").Copy()
[void]$ws.Range("A52:G52").PasteSpecial(-4122)

$ws.Cells.Item($newRow, 1).Value = "Linking_AutoUser"
$ws.Cells.Item($newRow, 2).Value = "Password1"
$ws.Cells.Item($newRow, 3).Value = ""
$ws.Cells.Item($newRow, 4).Value = ""
$ws.Cells.Item($newRow, 5).Value = "Default user for Linking tests"
$ws.Cells.Item($newRow, 6).Value = "N"
$ws.Cells.Item($newRow, 7).Value = "linking.autouser@mailinator.com"

# Restore the view state: sheet active, scrolled down to the new rows, with
# C22 selected (matches the author's cursor position when saving).
[void]$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
[void]$ws.Range("C22").Select()
